$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record row at row 145, pushing the existing rows
# 145..234 down to 146..235 (matches the dimension growing to A1:R235).
$ws.Rows("145:145").Insert()

# Populate the newly-inserted row with the new weekly record
# (same market/category/quality as the former row 145, new date + prices).
$ws.Range("A145").Value = 11
$ws.Range("B145").Value = "Vega Monumental Concepción"
$ws.Range("C145").Value = "Bíobío"
$ws.Range("D145").Value = 44673
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 100112008
$ws.Range("G145").Value = "Coliflor"
$ws.Range("H145").Value = "Sin especificar"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 2200
$ws.Range("K145").Value = 800
$ws.Range("L145").Value = 900
$ws.Range("M145").Value = 855
$ws.Range("N145").Value = "$/unidad"
$ws.Range("O145").Value = "Región Metropolitana"
$ws.Range("P145").Value = 855
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = "Hortaliza"
